$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells
$ws.Range("G1").Value = "liberal"
$ws.Range("H1").Value = "kapica"
$ws.Range("I1").Value = "peker"

# Match header formatting (bold) used by existing headers A1:F1
$ws.Range("G1:I1").Font.Bold = $true

# New data values in row 6
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 8

# Update selection to mirror the diff (selection moved to I7)
$ws.Range("I7").Select()
